$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures ---
$ws.Range("E11").Value2 = 414820    # VALOR MORA total
$ws.Range("F13").Value2 = 2         # Cant. Periodos

# --- Insert a new data row before the last (bottom-bordered) row so the
#     bottom-border row shifts from 19 -> 20, then fill all 5 data rows. ---
$ws.Rows("19:19").Insert()

# Give the newly inserted row 19 the same formatting as row 18 (the
# "middle" row style) before writing values into it.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)   # xlPasteFormats

# Row 16: ISABEL CRISTINA ALZATE MARIN, period 2507
$ws.Range("C16").Value2 = "1017182514"
$ws.Range("D16").Value2 = "ISABEL CRISTINA ALZATE MARIN"
$ws.Range("E16").Value2 = "2507"

# Row 17: RICK PETER HERNANDEZ RUSSO, period 2507
$ws.Range("C17").Value2 = "92642113"
$ws.Range("D17").Value2 = "RICK PETER HERNANDEZ RUSSO"
$ws.Range("E17").Value2 = "2507"
$ws.Range("F17").Value2 = 122000
$ws.Range("G17").Value2 = 3050000

# Row 18: JORGE ANDRES GARCIA OSORIO (new worker), period 2508
$ws.Range("C18").Value2 = "1001835668"
$ws.Range("D18").Value2 = "JORGE ANDRES GARCIA OSORIO"
$ws.Range("E18").Value2 = "2508"

# Row 19 (new): ISABEL CRISTINA ALZATE MARIN, period 2508
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1017182514"
$ws.Range("D19").Value2 = "ISABEL CRISTINA ALZATE MARIN"
$ws.Range("E19").Value2 = "2508"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500

# Row 20 (old row 19, shifted down): RICK PETER HERNANDEZ RUSSO, period 2508
$ws.Range("E20").Value2 = "2508"

Write-Host "Edits applied"
